# WorldID, ServerID -> 8-bit (uint8_t) change
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServerConfig")

# B3/C3 hold the type name for WorldID / ServerID columns ("int16_t" -> "uint8_t")
$ws.Range("B3").Value = "uint8_t"
$ws.Range("C3").Value = "uint8_t"

# Sample ServerID value shrinks from 1000 to 1 to fit the new 8-bit type
$ws.Range("C4").Value = 1

# Update the active selection to match the saved view state
$ws.Range("E10").Select()
